$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("LP1912")
$sheet1.Cells.Item(2,1).Value = "Última actualización: 17:13:12"
$sheet1.Cells.Item(3,1).Value = "Total filas: 381"
$sheet1.Cells.Item(52,1).Value = "07:36:59"
$sheet1.Cells.Item(52,3).Value = "17_ROMERO"
$sheet1.Cells.Item(52,4).Value = 26
$sheet1.Cells.Item(53,1).Value = "06:52:52"
$sheet1.Cells.Item(53,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(53,4).Value = 70
$sheet1.Cells.Item(64,1).Value = "06:45:50"
$sheet1.Cells.Item(64,3).Value = "14_ABASTO"
$sheet1.Cells.Item(64,4).Value = 104
$sheet1.Cells.Item(65,1).Value = "08:29:19"
$sheet1.Cells.Item(65,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(65,4).Value = 0
$sheet1.Cells.Item(87,3).Value = "11_ETCHEVERRY"
$sheet1.Cells.Item(88,3).Value = "15_ABASTO"
$sheet1.Cells.Item(149,1).Value = "11:11:31"
$sheet1.Cells.Item(149,3).Value = "16_SANTA ANA"
$sheet1.Cells.Item(149,4).Value = 0
$sheet1.Cells.Item(150,1).Value = "10:04:17"
$sheet1.Cells.Item(150,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(150,4).Value = 67
$sheet1.Cells.Item(151,1).Value = "10:36:18"
$sheet1.Cells.Item(151,3).Value = "15_ABASTO"
$sheet1.Cells.Item(151,4).Value = 35
$sheet1.Cells.Item(162,1).Value = "10:48:14"
$sheet1.Cells.Item(162,3).Value = "16_SANTA ANA"
$sheet1.Cells.Item(162,4).Value = 45
$sheet1.Cells.Item(163,1).Value = "10:36:18"
$sheet1.Cells.Item(163,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(163,4).Value = 57
$sheet1.Cells.Item(173,1).Value = "11:53:59"
$sheet1.Cells.Item(173,3).Value = "225_GOMEZ"
$sheet1.Cells.Item(173,4).Value = 0
$sheet1.Cells.Item(174,1).Value = "10:55:25"
$sheet1.Cells.Item(174,3).Value = "15_ABASTO"
$sheet1.Cells.Item(174,4).Value = 58
$sheet1.Cells.Item(183,3).Value = "16_SANTA ANA"
$sheet1.Cells.Item(184,3).Value = "15_ABASTO"
$sheet1.Cells.Item(194,1).Value = "10:36:18"
$sheet1.Cells.Item(194,3).Value = "11_ETCHEVERRY"
$sheet1.Cells.Item(194,4).Value = 114
$sheet1.Cells.Item(195,1).Value = "11:53:59"
$sheet1.Cells.Item(195,3).Value = "16_P MOR-SANTA ANA"
$sheet1.Cells.Item(195,4).Value = 37
$sheet1.Cells.Item(196,1).Value = "10:36:18"
$sheet1.Cells.Item(196,3).Value = "16_P MOR-SANTA ANA"
$sheet1.Cells.Item(196,4).Value = 115
$sheet1.Cells.Item(197,1).Value = "10:48:14"
$sheet1.Cells.Item(197,3).Value = "11_ETCHEVERRY"
$sheet1.Cells.Item(197,4).Value = 103
$sheet1.Cells.Item(220,1).Value = "11:11:31"
$sheet1.Cells.Item(220,3).Value = "14_ABASTO"
$sheet1.Cells.Item(220,4).Value = 115
$sheet1.Cells.Item(221,1).Value = "12:11:45"
$sheet1.Cells.Item(221,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(221,4).Value = 55
$sheet1.Cells.Item(232,3).Value = "16_SANTA ANA"
$sheet1.Cells.Item(233,3).Value = "17_ROMERO"
$sheet1.Cells.Item(244,1).Value = "13:39:24"
$sheet1.Cells.Item(244,3).Value = "16_SANTA ANA"
$sheet1.Cells.Item(244,4).Value = 0
$sheet1.Cells.Item(245,1).Value = "12:32:47"
$sheet1.Cells.Item(245,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(245,4).Value = 67
$sheet1.Cells.Item(246,1).Value = "11:53:59"
$sheet1.Cells.Item(246,3).Value = "17X38_ROMERO"
$sheet1.Cells.Item(246,4).Value = 106
$sheet1.Cells.Item(351,1).Value = "17:13:12"
$sheet1.Cells.Item(351,2).Value = "17:15"
$sheet1.Cells.Item(351,3).Value = "10_OLMOS"
$sheet1.Cells.Item(351,4).Value = 2
$sheet1.Cells.Item(352,2).Value = "17:17"
$sheet1.Cells.Item(352,3).Value = "17_ROMERO"
$sheet1.Cells.Item(352,4).Value = 92
$sheet1.Cells.Item(353,1).Value = "15:45:31"
$sheet1.Cells.Item(353,2).Value = "17:24"
$sheet1.Cells.Item(353,4).Value = 99
$sheet1.Cells.Item(354,1).Value = "16:43:14"
$sheet1.Cells.Item(354,2).Value = "17:25"
$sheet1.Cells.Item(354,3).Value = "11_ETCHEVERRY"
$sheet1.Cells.Item(354,4).Value = 42
$sheet1.Cells.Item(355,1).Value = "15:57:19"
$sheet1.Cells.Item(355,2).Value = "17:27"
$sheet1.Cells.Item(355,3).Value = "15_ABASTO"
$sheet1.Cells.Item(355,4).Value = 90
$sheet1.Cells.Item(356,1).Value = "16:36:34"
$sheet1.Cells.Item(356,2).Value = "17:30"
$sheet1.Cells.Item(356,4).Value = 54
$sheet1.Cells.Item(357,1).Value = "16:13:19"
$sheet1.Cells.Item(357,2).Value = "17:31"
$sheet1.Cells.Item(357,4).Value = 78
$sheet1.Cells.Item(358,1).Value = "16:52:31"
$sheet1.Cells.Item(358,2).Value = "17:33"
$sheet1.Cells.Item(358,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(358,4).Value = 41
$sheet1.Cells.Item(359,1).Value = "15:45:31"
$sheet1.Cells.Item(359,2).Value = "17:34"
$sheet1.Cells.Item(359,3).Value = "10_OLMOS"
$sheet1.Cells.Item(359,4).Value = 109
$sheet1.Cells.Item(360,1).Value = "16:13:19"
$sheet1.Cells.Item(360,3).Value = "27_EL RETIRO"
$sheet1.Cells.Item(360,4).Value = 82
$sheet1.Cells.Item(361,2).Value = "17:35"
$sheet1.Cells.Item(361,3).Value = "16_P MOR-SANTA ANA"
$sheet1.Cells.Item(361,4).Value = 110
$sheet1.Cells.Item(362,1).Value = "15:45:31"
$sheet1.Cells.Item(362,2).Value = "17:36"
$sheet1.Cells.Item(362,4).Value = 111
$sheet1.Cells.Item(363,1).Value = "16:43:14"
$sheet1.Cells.Item(363,2).Value = "17:37"
$sheet1.Cells.Item(363,3).Value = "27_EL RETIRO"
$sheet1.Cells.Item(363,4).Value = 54
$sheet1.Cells.Item(364,1).Value = "15:45:31"
$sheet1.Cells.Item(364,2).Value = "17:38"
$sheet1.Cells.Item(364,3).Value = "17X38_ROMERO"
$sheet1.Cells.Item(364,4).Value = 113
$sheet1.Cells.Item(365,1).Value = "16:27:37"
$sheet1.Cells.Item(365,2).Value = "17:44"
$sheet1.Cells.Item(365,4).Value = 77
$sheet1.Cells.Item(366,1).Value = "15:57:19"
$sheet1.Cells.Item(366,2).Value = "17:45"
$sheet1.Cells.Item(366,3).Value = "215B_EL PATO"
$sheet1.Cells.Item(366,4).Value = 108
$sheet1.Cells.Item(367,1).Value = "16:43:14"
$sheet1.Cells.Item(367,2).Value = "17:47"
$sheet1.Cells.Item(367,3).Value = "16_SANTA ANA"
$sheet1.Cells.Item(367,4).Value = 64
$sheet1.Cells.Item(368,1).Value = "16:27:37"
$sheet1.Cells.Item(368,2).Value = "17:49"
$sheet1.Cells.Item(368,3).Value = "17X38_ROMERO"
$sheet1.Cells.Item(368,4).Value = 82
$sheet1.Cells.Item(369,1).Value = "15:57:19"
$sheet1.Cells.Item(369,2).Value = "17:51"
$sheet1.Cells.Item(369,3).Value = "215_EL PELIGRO"
$sheet1.Cells.Item(369,4).Value = 114
$sheet1.Cells.Item(370,1).Value = "17:13:12"
$sheet1.Cells.Item(370,2).Value = "17:59"
$sheet1.Cells.Item(370,3).Value = "16_SANTA ANA"
$sheet1.Cells.Item(370,4).Value = 46
$sheet1.Cells.Item(371,1).Value = "16:27:37"
$sheet1.Cells.Item(371,2).Value = "18:02"
$sheet1.Cells.Item(371,3).Value = "17_ROMERO"
$sheet1.Cells.Item(371,4).Value = 95
$sheet1.Cells.Item(372,2).Value = "18:03"
$sheet1.Cells.Item(372,3).Value = "17_ROMERO"
$sheet1.Cells.Item(372,4).Value = 110
$sheet1.Cells.Item(373,1).Value = "16:36:34"
$sheet1.Cells.Item(373,2).Value = "18:04"
$sheet1.Cells.Item(373,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(373,4).Value = 88
$sheet1.Cells.Item(374,1).Value = "16:13:19"
$sheet1.Cells.Item(374,2).Value = "18:04"
$sheet1.Cells.Item(374,3).Value = "14_ABASTO"
$sheet1.Cells.Item(374,4).Value = 111
$sheet1.Cells.Item(375,2).Value = "18:05"
$sheet1.Cells.Item(375,3).Value = "14_ABASTO"
$sheet1.Cells.Item(375,4).Value = 82
$sheet1.Cells.Item(376,1).Value = "17:13:12"
$sheet1.Cells.Item(376,2).Value = "18:12"
$sheet1.Cells.Item(376,3).Value = "23_HERNANDEZ"
$sheet1.Cells.Item(376,4).Value = 59
$sheet1.Cells.Item(377,1).Value = "17:13:12"
$sheet1.Cells.Item(377,2).Value = "18:14"
$sheet1.Cells.Item(377,3).Value = "10_OLMOS"
$sheet1.Cells.Item(377,4).Value = 61
$sheet1.Cells.Item(378,1).Value = "16:27:37"
$sheet1.Cells.Item(378,2).Value = "18:24"
$sheet1.Cells.Item(378,3).Value = "11_ETCHEVERRY"
$sheet1.Cells.Item(378,4).Value = 117
$sheet1.Cells.Item(379,1).Value = "16:43:14"
$sheet1.Cells.Item(379,2).Value = "18:25"
$sheet1.Cells.Item(379,3).Value = "11_ETCHEVERRY"
$sheet1.Cells.Item(379,4).Value = 102
$sheet1.Cells.Item(379,5).Value = "LP1912"
$sheet1.Cells.Item(380,1).Value = "17:13:12"
$sheet1.Cells.Item(380,2).Value = "18:27"
$sheet1.Cells.Item(380,3).Value = "15_ABASTO"
$sheet1.Cells.Item(380,4).Value = 74
$sheet1.Cells.Item(380,5).Value = "LP1912"
$sheet1.Cells.Item(381,1).Value = "16:36:34"
$sheet1.Cells.Item(381,2).Value = "18:34"
$sheet1.Cells.Item(381,3).Value = "14X44_ABASTO"
$sheet1.Cells.Item(381,4).Value = 118
$sheet1.Cells.Item(381,5).Value = "LP1912"
$sheet1.Cells.Item(382,1).Value = "16:43:14"
$sheet1.Cells.Item(382,2).Value = "18:38"
$sheet1.Cells.Item(382,3).Value = "17X38_ROMERO"
$sheet1.Cells.Item(382,4).Value = 115
$sheet1.Cells.Item(382,5).Value = "LP1912"
$sheet1.Cells.Item(383,1).Value = "16:43:14"
$sheet1.Cells.Item(383,2).Value = "18:41"
$sheet1.Cells.Item(383,3).Value = "16_P MOR-SANTA ANA"
$sheet1.Cells.Item(383,4).Value = 118
$sheet1.Cells.Item(383,5).Value = "LP1912"
$sheet1.Cells.Item(384,1).Value = "17:13:12"
$sheet1.Cells.Item(384,2).Value = "18:41"
$sheet1.Cells.Item(384,3).Value = "14_ABASTO"
$sheet1.Cells.Item(384,4).Value = 88
$sheet1.Cells.Item(384,5).Value = "LP1912"
$sheet1.Cells.Item(385,1).Value = "17:13:12"
$sheet1.Cells.Item(385,2).Value = "19:01"
$sheet1.Cells.Item(385,3).Value = "17_ROMERO"
$sheet1.Cells.Item(385,4).Value = 108
$sheet1.Cells.Item(385,5).Value = "LP1912"
$sheet1.Cells.Item(386,1).Value = "17:13:12"
$sheet1.Cells.Item(386,2).Value = "19:11"
$sheet1.Cells.Item(386,3).Value = "81_EL PELIGRO"
$sheet1.Cells.Item(386,4).Value = 118
$sheet1.Cells.Item(386,5).Value = "LP1912"

$sheet2 = $wb.Worksheets.Item("LP1912-215")
$sheet2.Cells.Item(2,1).Value = "Última actualización: 17:13:12"

$sheet3 = $wb.Worksheets.Item("6203-6173")
$sheet3.Cells.Item(2,1).Value = "Última actualización: 17:13:12"
$sheet3.Cells.Item(3,1).Value = "Total filas: 45"
$sheet3.Cells.Item(49,1).Value = "17:13:12"
$sheet3.Cells.Item(49,2).Value = "18:26"
$sheet3.Cells.Item(49,3).Value = "215C_LA PLATA"
$sheet3.Cells.Item(49,4).Value = 73
$sheet3.Cells.Item(49,5).Value = "L6203"
$sheet3.Cells.Item(50,1).Value = "17:13:12"
$sheet3.Cells.Item(50,2).Value = "19:11"
$sheet3.Cells.Item(50,3).Value = "215B_LP-P MOR-1 Y 57"
$sheet3.Cells.Item(50,4).Value = 118
$sheet3.Cells.Item(50,5).Value = "L6173"

